$d = $word.ActiveDocument

# --- 1) Insert a "_GoBack" bookmark right after the
#        "(Determining the maximum No. assignments possible)" run. ---
$rng = $d.Content
$found = $rng.Find.Execute(
    "(Determining the maximum No. assignments possible)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPoint = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $insertPoint)
}

# --- 2) Fix the spacing typo "유지 하면서" -> "유지하면서" in the dual-feasible line. ---
$null = $d.Content.Find.Execute(
    "유지 하면서",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "유지하면서", 2)

# --- 3) Merge the run split around the old "_GoBack" bookmark
#        ("Machine의 수 일 때, make dummy rows or c" | bookmark | "olumns such that ")
#        back into a single run, which also removes that now-stale bookmark. ---
$null = $d.Content.Find.Execute(
    "Machine의 수 일 때, make dummy rows or columns such that ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Machine의 수 일 때, make dummy rows or columns such that ", 2)
